# issue #5: add legislator_id, name, date into dataframe
#
# The 股票 (stock) sheet gains three new trailing columns:
#   H: date             -> "2012-04-09"
#   I: legislator_name  -> "林佳龍"
#   J: legislator_id    -> 1741
# for every existing data row (2..11), plus matching headers in row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)   # 股票

$lastRow = 11                 # existing data rows are 2..11

# --- Header row: copy the look of the existing header cells (B1:G1) onto H1:J1
$ws.Range("B1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# --- Data rows: copy the look of an existing plain data cell (B2) onto the
#     new cells so the new columns match the rest of the row (not column A's
#     bold/bordered style).
$ws.Range("B2").Copy()
$ws.Range("H2:J$lastRow").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

for ($r = 2; $r -le $lastRow; $r++) {
    # Force text so "2012-04-09" isn't auto-converted into a date serial.
    $ws.Cells.Item($r, 8).NumberFormat = "@"
    $ws.Cells.Item($r, 8).Value = "2012-04-09"        # H -> date
    $ws.Cells.Item($r, 9).Value = "林佳龍"              # I -> legislator_name
    $ws.Cells.Item($r, 10).Value = 1741                # J -> legislator_id
}

# Re-apply the plain data-cell look to the date column so the temporary text
# number-format doesn't leave a distinct style behind.
$ws.Range("B2").Copy()
$ws.Range("H2:H$lastRow").PasteSpecial(-4122)
$excel.CutCopyMode = 0
